$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 366.8974733393968
$ws.Range("C2").Value = 4951.183334891753
$ws.Range("D2").Value = 3472.370637841082
